# Updates Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.240.83'
$ws.Cells.Item(2, 5).Value = '  +1.29%  '
$ws.Cells.Item(3, 4).Value = '1.900.33'
$ws.Cells.Item(3, 5).Value = '  +0.74%  '
$ws.Cells.Item(4, 4).Value = '0.9987'
$ws.Cells.Item(4, 5).Value = '  -0.27%  '
$ws.Cells.Item(5, 4).Value = '0.7482'
$ws.Cells.Item(5, 5).Value = '  +0.11%  '
$ws.Cells.Item(6, 4).Value = '242.86'
$ws.Cells.Item(6, 5).Value = '  +0.25%  '
$ws.Cells.Item(7, 4).Value = '0.9974'
$ws.Cells.Item(7, 5).Value = '  -0.41%  '
$ws.Cells.Item(8, 4).Value = '0.3180'
$ws.Cells.Item(8, 5).Value = '  +1.88%  '
$ws.Cells.Item(9, 4).Value = '0.07269'
$ws.Cells.Item(9, 5).Value = '  +2.09%  '
$ws.Cells.Item(10, 4).Value = '25.17'
$ws.Cells.Item(10, 5).Value = '  -0.98%  '
$ws.Cells.Item(11, 4).Value = '0.08378'
$ws.Cells.Item(11, 5).Value = '  -1.35%  '
$ws.Cells.Item(12, 4).Value = '0.7652'
$ws.Cells.Item(12, 5).Value = '  +0.76%  '
$ws.Cells.Item(13, 4).Value = '5.469'
$ws.Cells.Item(13, 5).Value = '  +2.17%  '
$ws.Cells.Item(14, 4).Value = '1.889.94'
$ws.Cells.Item(14, 5).Value = '  -0.67%  '
$ws.Cells.Item(15, 4).Value = '93.27'
$ws.Cells.Item(15, 5).Value = '  -0.16%  '
$ws.Cells.Item(16, 4).Value = '6.196'
$ws.Cells.Item(16, 5).Value = '  +0.78%  '
$ws.Cells.Item(17, 4).Value = '30.195.32'
$ws.Cells.Item(17, 5).Value = '  +0.78%  '
$ws.Cells.Item(18, 4).Value = '251.78'
$ws.Cells.Item(18, 5).Value = '  +3.58%  '
$ws.Cells.Item(19, 4).Value = '13.72'
$ws.Cells.Item(19, 5).Value = '  +0.01%  '
$ws.Cells.Item(20, 4).Value = '0.000007894'
$ws.Cells.Item(20, 5).Value = '  +1.34%  '
$ws.Cells.Item(21, 4).Value = '2.167.45'
$ws.Cells.Item(21, 5).Value = '  -0.55%  '
$ws.Cells.Item(22, 4).Value = '0.9971'
$ws.Cells.Item(22, 5).Value = '  -0.28%  '
$ws.Cells.Item(23, 4).Value = '8.046'
$ws.Cells.Item(23, 5).Value = '  +0.36%  '
$ws.Cells.Item(24, 4).Value = '0.9985'
$ws.Cells.Item(24, 5).Value = '  -0.30%  '
$ws.Cells.Item(25, 4).Value = '0.1600'
$ws.Cells.Item(25, 5).Value = '  +0.48%  '
$ws.Cells.Item(26, 4).Value = '9.344'
$ws.Cells.Item(26, 5).Value = '  -0.33%  '
$ws.Cells.Item(27, 4).Value = '164.72'
$ws.Cells.Item(27, 5).Value = '  +1.10%  '
$ws.Cells.Item(28, 4).Value = '18.87'
$ws.Cells.Item(28, 5).Value = '  +0.61%  '
$ws.Cells.Item(29, 4).Value = '2.077'
$ws.Cells.Item(29, 5).Value = '  +2.45%  '
$ws.Cells.Item(30, 4).Value = '1.479'
$ws.Cells.Item(30, 5).Value = '  -2.19%  '
$ws.Cells.Item(31, 4).Value = '4.632'
$ws.Cells.Item(31, 5).Value = '  +3.56%  '
$ws.Cells.Item(32, 4).Value = '1.539'
$ws.Cells.Item(32, 5).Value = '  +0.68%  '
$ws.Cells.Item(33, 4).Value = '4.243'
$ws.Cells.Item(33, 5).Value = '  +3.51%  '
$ws.Cells.Item(34, 4).Value = '0.05421'
$ws.Cells.Item(34, 5).Value = '  +0.46%  '
$ws.Cells.Item(35, 4).Value = '1.263'
$ws.Cells.Item(35, 5).Value = '  +2.14%  '
$ws.Cells.Item(36, 4).Value = '0.7662'
$ws.Cells.Item(36, 5).Value = '  +3.06%  '
$ws.Cells.Item(37, 4).Value = '0.9960'
$ws.Cells.Item(37, 5).Value = '  -0.74%  '
$ws.Cells.Item(38, 4).Value = '2.713'
$ws.Cells.Item(38, 5).Value = '  -0.02%  '
$ws.Cells.Item(39, 4).Value = '0.01979'
$ws.Cells.Item(39, 5).Value = '  +2.51%  '
$ws.Cells.Item(40, 4).Value = '2.774'
$ws.Cells.Item(40, 5).Value = '  +0.19%  '
$ws.Cells.Item(41, 4).Value = '0.4603'
$ws.Cells.Item(41, 5).Value = '  +3.35%  '
$ws.Cells.Item(42, 4).Value = '1.099.79'
$ws.Cells.Item(42, 5).Value = '  +0.49%  '
$ws.Cells.Item(43, 4).Value = '6.096'
$ws.Cells.Item(43, 5).Value = '  +0.34%  '
$ws.Cells.Item(44, 4).Value = '73.02'
$ws.Cells.Item(44, 5).Value = '  +0.69%  '
$ws.Cells.Item(45, 4).Value = '0.8740'
$ws.Cells.Item(45, 5).Value = '  +2.10%  '
$ws.Cells.Item(46, 4).Value = '104.59'
$ws.Cells.Item(46, 5).Value = '  +2.11%  '
$ws.Cells.Item(47, 4).Value = '0.9995'
$ws.Cells.Item(47, 5).Value = '  -0.15%  '
$ws.Cells.Item(48, 4).Value = '1.879'
$ws.Cells.Item(48, 5).Value = '  +0.97%  '
$ws.Cells.Item(49, 4).Value = '7.652'
$ws.Cells.Item(49, 5).Value = '  -0.04%  '
$ws.Cells.Item(50, 4).Value = '9.644'
$ws.Cells.Item(50, 5).Value = '  -0.83%  '
$ws.Cells.Item(51, 4).Value = '2.047.54'
$ws.Cells.Item(51, 5).Value = '  -0.13%  '
